# 建筑业企业签订合同和承包工程完成情况.xlsx
# Remove the oldest six years of data (2004年-2009年, rows 2-7) so that the
# table now starts at 2010年, then append the two new years (2021年 and
# 2022年) that were added at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 2..7 (2004年..2009年); this shifts 2010年..2020年 up to rows 2..12.
$ws.Range("A2:A7").EntireRow.Delete()

# Row 12 (2020年) keeps the header-style formatting (bold + border) used by
# column A; copy it down onto the two new rows so the new year labels match
# the look of the existing ones.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A12").Copy($ws.Range("A14"))

# New row 13: 2021年
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 3123272991.8
$ws.Range("C13").Value = 194632661.3
$ws.Range("D13").Value = 123451527.4
$ws.Range("E13").Value = 3445591950.9
$ws.Range("F13").Value = 2859601884.7
$ws.Range("G13").Value = 6568864942.7
$ws.Range("H13").Value = 2736150357.3

# New row 14: 2022年 (only partial data is available for this year)
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 3491933383.6
$ws.Range("E14").Value = 3664813494.7
$ws.Range("G14").Value = 7156746878.3
